$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.198.26"
$ws.Range("E2").Value = "  -3.19%  "

$ws.Range("D3").Value = "1.607.97"
$ws.Range("E3").Value = "  -2.68%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "303.00"
$ws.Range("E6").Value = "  -2.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3767"
$ws.Range("E7").Value = "  -3.33%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3651"
$ws.Range("E8").Value = "  -4.71%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.67"
$ws.Range("E9").Value = "  -5.14%  "

$ws.Range("E10").Value = "  -0.02%  "

$ws.Range("E11").Value = "  -6.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08082"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.02"
$ws.Range("E13").Value = "  -3.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.582"
$ws.Range("E14").Value = "  -7.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.639"
$ws.Range("E15").Value = "  -3.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001269"

$ws.Range("D17").Value = "1.609.11"
$ws.Range("E17").Value = "  -2.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.48"
$ws.Range("E18").Value = "  -3.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06786"
$ws.Range("E19").Value = "  -2.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.36"
$ws.Range("E20").Value = "  -6.89%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.580"
$ws.Range("E21").Value = "  -5.13%  "

$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.11"
$ws.Range("E23").Value = "  -4.55%  "

$ws.Range("D24").Value = "23.223.26"
$ws.Range("E24").Value = "  -3.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.353"
$ws.Range("E25").Value = "  -4.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.908"
$ws.Range("E26").Value = "  -2.26%  "

$ws.Range("E27").Value = "  -4.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.32"
$ws.Range("E28").Value = "  -0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.267"
$ws.Range("E29").Value = "  -3.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.47"
$ws.Range("E30").Value = "  -4.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.419"
$ws.Range("E31").Value = "  -2.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.925"
$ws.Range("E32").Value = "  -11.53%  "

$ws.Range("D33").Value = "1.785.64"
$ws.Range("E33").Value = "  -2.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9782"
$ws.Range("E34").Value = "  -6.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07729"
$ws.Range("E35").Value = "  -4.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02785"
$ws.Range("E36").Value = "  -5.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.269"
$ws.Range("E37").Value = "  -7.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2551"
$ws.Range("E38").Value = "  -5.05%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.10"
$ws.Range("E39").Value = "  -7.34%  "

$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08859"
$ws.Range("E40").Value = "  -3.06%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.396"
$ws.Range("E41").Value = "  -1.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7160"
$ws.Range("E42").Value = "  -5.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.78"
$ws.Range("E43").Value = "  -5.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.87"
$ws.Range("E44").Value = "  -2.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6605"
$ws.Range("E45").Value = "  -4.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("E47").Value = "  -6.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.984"
$ws.Range("E48").Value = "  -2.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08012"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.33"
$ws.Range("E50").Value = "  -2.59%  "

$ws.Range("E51").Value = "  -3.36%  "
